$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "input" column (F) header - copy formatting from the neighboring
# "last_run" header (E1) so it matches the other bold/bordered headers.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "input"

# Update existing "output" column (D) values
$ws.Range("D2").Value = "PDF"
$ws.Range("D3").Value = "EXCEL"
$ws.Range("D4").Value = "EXCEL"

# Fill new "input" column (F) values
$ws.Range("F2").Value = "EXCEL"
$ws.Range("F3").Value = "PDF"
$ws.Range("F4").Value = "EXCEL"

$excel.CutCopyMode = 0
